# Edit TableS1_Gsp1_HRas_annotations.xlsx
# 1. Merge the "GTPase regions" and "Contacts Nucleotide" rows into a single
#    "Active site regions" row (delete the "Contacts Nucleotide" row, and
#    update the remaining row's label + activating-HRas-positions cell).
# 2. Set the print page orientation to portrait.
# 3. Update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3 ("GTPase regions" -> "Active site regions") ---
$ws.Range("A3").Value2 = "Active site regions"
$ws.Range("C3").Value2 = "G12, G13, V14, K16, A18, P34, T58, A59, G60, Q61, E63, R68, N116, K117, D119, L120, S145, A146, K147"

# --- Delete row 4 ("Contacts Nucleotide"), shifting rows 5-7 up ---
$ws.Rows.Item(4).Delete()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1 | Out-Null

# --- Update the worksheet's active cell selection ---
$ws.Range("H10").Select() | Out-Null
